# "Main sequence -> Main flowchart"
# Adds a new "downloadPath" setting row to the Settings sheet and makes
# Settings the active/selected sheet (it was previously "Assets").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Bring Settings to the front / make it the active sheet (this also clears
# the previously-active "Assets" sheet's tabSelected flag).
$ws.Activate()

# New row of settings data.
$ws.Range("A2").Value = "downloadPath"
$ws.Range("B2").Value = "Data\Downloads"

# Widen columns A and B to fit the new content (column C already had a
# fitted width).
$ws.Columns("A:B").AutoFit() | Out-Null

# Leave the same kind of "last used cell" selection state recorded by Excel.
$ws.Range("O15").Select() | Out-Null
